# Apply the committed edits to Exp_data sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 185: update forecast-ish values for years 2020-2022 (AR/AS/AT) ---
$ws.Range("AR185").Value = 12
$ws.Range("AS185").Value = 13.3
$ws.Range("AT185").Value = 12.2

# --- Rows 189-191: add a new year-2024 column (AV) value ---
$ws.Range("AV189").Value = 4533
$ws.Range("AV190").Value = 1258
$ws.Range("AV191").Value = 1122

# --- Remove the four "Z16_B01 crime-offence" detail rows (old 214-217) ---
# This shifts every row below them up by four (old 218 -> new 214, ... old 229 -> new 225).
$ws.Rows("214:217").Delete()

# --- Row 213: replace the "Z16_B01...D001" record with a new "Z15_B04...D001" record ---
$ws.Range("A213").Value = "Z15_B04_P01_Ib01_I01_D001"
$ws.Range("B213").Value = "Z15_B04_P01_Ib01_I01"
$ws.Range("I213").Value = "E_HA"
$ws.Range("L213").Value = $false

# Drop the old yearly time series (1993-2020) entirely - the new record has no history
$ws.Range("Q213:AR213").ClearContents()

# New record only carries placeholder values for 2021-2023
$ws.Range("AS213").Value = 999
$ws.Range("AT213").Value = 999
$ws.Range("AU213").Value = 999

$ws.Range("CP213").Value = "Z15_B04_P01_Ib01"
